# Update countries & provincias Spain
# Applies the data refresh described in the commit/diff:
#  - Reorders "Japon" to sit before "Peru" (Chequia, Japon, Peru, Rumania, Ecuador)
#  - Reorders "Consejo Danes para los Refugiados" to sit before "Guinea"
#  - Updates the "Datos actualizados" timestamp
#  - Refreshes several countries' case statistics

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Timestamp (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 10 de Abril de 2020 a las 05:52"

# --- Australia (row 27) : isolated stat tweak ---
$ws.Range("E27").Value = 3112
$ws.Range("G27").Value = 2
$ws.Range("H27").Value = 53

# --- Reorder: Chequia(31) / Japon(32) / Peru(33) / Rumania(34) / Ecuador(35) ---
$ws.Range("A32").Value = "Japon"
$ws.Range("B32").Value = 5530
$ws.Range("C32").Value = 183
$ws.Range("D32").Value = 685
$ws.Range("E32").Value = 4746
$ws.Range("F32").Value = 109
$ws.Range("G32").Value = 0
$ws.Range("H32").Value = 99

$ws.Range("A33").Value = "Peru"
$ws.Range("B33").Value = 5256
$ws.Range("C33").Value = 0
$ws.Range("D33").Value = 1438
$ws.Range("E33").Value = 3680
$ws.Range("F33").Value = 124
$ws.Range("G33").Value = 0
$ws.Range("H33").Value = 138

$ws.Range("A34").Value = "Rumania"
$ws.Range("B34").Value = 5202
$ws.Range("C34").Value = 0
$ws.Range("D34").Value = 647
$ws.Range("E34").Value = 4307
$ws.Range("F34").Value = 178
$ws.Range("G34").Value = 0
$ws.Range("H34").Value = 248

# Row 35 (Ecuador) is unchanged.

# --- Nueva Zelanda (row 63) : isolated stat tweak ---
$ws.Range("F63").Value = 5

# --- Reorder: Georgia(112) / Consejo Danes(113) / Guinea(114) / Sri Lanka(115) /
#              Isla de Man(116) / Kenia(117) / Mayotte(118) / Islas Feroe(119) / Venezuela(120) ---
$ws.Range("A113").Value = "Consejo Danes para los Refugiados"
$ws.Range("B113").Value = 215
$ws.Range("C113").Value = 35
$ws.Range("D113").Value = 13
$ws.Range("E113").Value = 182
$ws.Range("F113").Value = 0
$ws.Range("G113").Value = 2
$ws.Range("H113").Value = 20

$ws.Range("A114").Value = "Guinea"
$ws.Range("B114").Value = 194
$ws.Range("C114").Value = 0
$ws.Range("D114").Value = 11
$ws.Range("E114").Value = 183
$ws.Range("F114").Value = 0
$ws.Range("G114").Value = 0
$ws.Range("H114").Value = 0

$ws.Range("A115").Value = "Sri Lanka"
$ws.Range("B115").Value = 190
$ws.Range("C115").Value = 0
$ws.Range("D115").Value = 49
$ws.Range("E115").Value = 134
$ws.Range("F115").Value = 5
$ws.Range("G115").Value = 0
$ws.Range("H115").Value = 7

$ws.Range("A116").Value = "Isla de Man"
$ws.Range("B116").Value = 190
$ws.Range("C116").Value = 0
$ws.Range("D116").Value = 92
$ws.Range("E116").Value = 97
$ws.Range("F116").Value = 12
$ws.Range("G116").Value = 0
$ws.Range("H116").Value = 1

$ws.Range("A117").Value = "Kenia"
$ws.Range("B117").Value = 184
$ws.Range("C117").Value = 0
$ws.Range("D117").Value = 12
$ws.Range("E117").Value = 165
$ws.Range("F117").Value = 2
$ws.Range("G117").Value = 0
$ws.Range("H117").Value = 7

$ws.Range("A118").Value = "Mayotte"
$ws.Range("B118").Value = 184
$ws.Range("C118").Value = 0
$ws.Range("D118").Value = 26
$ws.Range("E118").Value = 156
$ws.Range("F118").Value = 4
$ws.Range("G118").Value = 0
$ws.Range("H118").Value = 2

$ws.Range("A119").Value = "Islas Feroe"
$ws.Range("B119").Value = 184
$ws.Range("C119").Value = 0
$ws.Range("D119").Value = 136
$ws.Range("E119").Value = 48
$ws.Range("F119").Value = 0
$ws.Range("G119").Value = 0
$ws.Range("H119").Value = 0

# Row 120 (Venezuela) is unchanged.
